$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 37)
$ws.Range("D2").Value = [double]"8.811544254521116E-08"
$ws.Range("E2").Value = [double]"8.811544254521116E-08"

# Row 3 (Control 4)
$ws.Range("D3").Value = [double]"3.568688884333009E-08"
$ws.Range("E3").Value = [double]"3.568688884333009E-08"

# Row 4 (Control 45)
$ws.Range("D4").Value = [double]"1.31491455725481E-34"
$ws.Range("E4").Value = [double]"1.31491455725481E-34"

# Row 5 (Control 48)
$ws.Range("D5").Value = [double]"1.110744401486915E-56"
$ws.Range("E5").Value = [double]"1.110744401486915E-56"

# Row 6 (Control 20)
$ws.Range("D6").Value = [double]"1.183103460400268E-10"
$ws.Range("E6").Value = [double]"1.183103460400268E-10"

# Row 7 (MDD 37)
$ws.Range("D7").Value = [double]"0.9999999999954301"
$ws.Range("E7").Value = [double]"4.569900013962069E-12"

# Row 8 (MDD 24)
$ws.Range("D8").Value = [double]"1.777229895013519E-12"
$ws.Range("E8").Value = [double]"0.9999999999982228"

# Row 9 (MDD 6)
$ws.Range("D9").Value = [double]"3.805269517419362E-11"
$ws.Range("E9").Value = [double]"0.9999999999619473"

# Row 10 (MDD 54)
$ws.Range("D10").Value = [double]"0.9999986531823306"
$ws.Range("E10").Value = [double]"1.346817669389822E-06"

# Row 11 (MDD 21)
$ws.Range("D11").Value = [double]"3.141077440158784E-10"
$ws.Range("E11").Value = [double]"0.9999999996858923"
$ws.Range("F11").Value = [double]"7.292929649353027"
